# Fixed naive component forecaster bug - Presentation state 11.02.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the stray C2 value (bug: this forecast cell should not have been populated)
$ws.Range("C2").ClearContents()

# Update recalculated forecast values (tiny floating point corrections from the
# fixed naive-component forecaster logic)
$ws.Range("C3").Value = -5.478010998490146
$ws.Range("E4").Value = -1.305195642355672
$ws.Range("C5").Value = 3.371423250978833
$ws.Range("E5").Value = 0.806063216063202
$ws.Range("E6").Value = 2.76635821344573
$ws.Range("C7").Value = -0.03183655677960751
$ws.Range("E7").Value = 1.102200073559856
$ws.Range("C8").Value = 1.812248956008777
$ws.Range("E8").Value = 1.209672013646301
$ws.Range("C10").Value = 1.575690123464613
$ws.Range("E10").Value = 1.643656926428538
$ws.Range("C12").Value = 2.337818484846443
$ws.Range("C13").Value = 0.8311911554373275
$ws.Range("E13").Value = 1.758956425699276
$ws.Range("C14").Value = -1.538034740964356
$ws.Range("E14").Value = -0.7351085756681197
$ws.Range("C15").Value = -0.5490727792360039
$ws.Range("C16").Value = 1.687572871803722
$ws.Range("C17").Value = -0.02880469535951891
$ws.Range("C18").Value = -0.1645795020818963
$ws.Range("E18").Value = 0.2574142441027716
$ws.Range("C19").Value = 0.06409464788890151
$ws.Range("E19").Value = -0.06354501920062816
